# Update report_co_so workbook:
#  - Sheet "CHI TIẾT DOANH THU": reorder A/B/C columns -> A: Ngày thực hiện, B: Tiền tố, C: Mã dịch vụ
#  - Insert new sheet "CHI TIẾT CHI TIÊU" (detailed expense log) right after "CHI TIẾT DOANH THU"
#  - Rename "DAONH SỐ CÁ NHÂN" -> "DOANH SỐ CÁ NHÂN"
#  - Rename "CHI TIÊU" -> "CHI TIÊU TỔNG HỢP"
#  - "LŨY KẾ NGÀY" stays the same (only shifts position)
#
# The remaining (unchanged-content) sheets are dropped and rebuilt in the
# correct order so that their internal sheetId numbering stays sequential,
# matching what a fresh re-generation of the report would produce.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Forces a value that Excel would otherwise auto-parse (e.g. a
    # dd-mm-yyyy-looking string) to be kept as plain text.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

function Write-Row($ws, [int]$rowIndex, [object[]]$values) {
    for ($col = 1; $col -le $values.Length; $col++) {
        $val = $values[$col - 1]
        if ($val -eq $null) {
            continue
        }
        $cell = $ws.Cells.Item($rowIndex, $col)
        if ($val -is [string] -and $val -match '^\d{2}-\d{2}-\d{4}$') {
            Set-TextValue $cell $val
        } else {
            $cell.Value = $val
        }
    }
}

function Set-StandardMargins($ws) {
    # Match the margins used throughout the rest of the workbook
    # (0.75in / 0.75in / 1in / 1in / 0.5in / 0.5in), expressed in points.
    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36
}

# ---------------------------------------------------------------------------
# 1. Fix column order on "CHI TIẾT DOANH THU" (first sheet): swap A/B/C so
#    that the execution date leads, followed by prefix, followed by code.
# ---------------------------------------------------------------------------
$wsDoanhThu = $wb.Worksheets.Item(1)

# Row 1 (header)
Set-TextValue $wsDoanhThu.Range("A1") "Ngày thực hiện"
Set-TextValue $wsDoanhThu.Range("B1") "Tiền tố"
Set-TextValue $wsDoanhThu.Range("C1") "Mã dịch vụ"

# Row 2
Set-TextValue $wsDoanhThu.Range("A2") "07-01-2024"
Set-TextValue $wsDoanhThu.Range("B2") "HD-LUXURY"
$wsDoanhThu.Range("C2").Value = 507

# Row 3
Set-TextValue $wsDoanhThu.Range("A3") "07-02-2024"
Set-TextValue $wsDoanhThu.Range("B3") "HD-LUXURY"
$wsDoanhThu.Range("C3").Value = 511

# ---------------------------------------------------------------------------
# 2. Drop the other 3 sheets (their data is unaffected, only name/position
#    change) so that they - and the new sheet - can be rebuilt in final
#    order with clean, sequential internal sheet ids.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(4).Delete() | Out-Null   # LŨY KẾ NGÀY
$wb.Worksheets.Item(3).Delete() | Out-Null   # CHI TIÊU
$wb.Worksheets.Item(2).Delete() | Out-Null   # DAONH SỐ CÁ NHÂN

# ---------------------------------------------------------------------------
# 3. Re-create sheets in the desired final order.
# ---------------------------------------------------------------------------

# 3a. "CHI TIẾT CHI TIÊU" (brand new detailed expense sheet)
$wsChiTietChiTieu = $wb.Worksheets.Add($null, $wsDoanhThu)
$wsChiTietChiTieu.Name = "CHI TIẾT CHI TIÊU"
Set-StandardMargins $wsChiTietChiTieu

Write-Row $wsChiTietChiTieu 1 @("Tiền tố", "Mã chi tiêu", "Ngày chi", "Cơ sở", "Phân loại", "Lượng chi")
Write-Row $wsChiTietChiTieu 2  @("CT", 571, "07-01-2024", "LONG XUYÊN", "Ứng Lương", 5000000)
Write-Row $wsChiTietChiTieu 3  @("CT", 572, "07-01-2024", "LONG XUYÊN", "Chi Phí CTV", 1800000)
Write-Row $wsChiTietChiTieu 4  @("CT", 574, "07-02-2024", "LONG XUYÊN", "Chi Phí Sinh Hoạt Tại Cơ Sở", 400000)
Write-Row $wsChiTietChiTieu 5  @("CT", 575, "07-02-2024", "LONG XUYÊN", "Chi Phí Sinh Hoạt Tại Cơ Sở", 100000)
Write-Row $wsChiTietChiTieu 6  @("CT", 576, "07-02-2024", "LONG XUYÊN", "Chi Phí Sinh Hoạt Tại Cơ Sở", 30000)
Write-Row $wsChiTietChiTieu 7  @("CT", 583, "07-03-2024", "LONG XUYÊN", "Chi Phí CTV", 1280000)
Write-Row $wsChiTietChiTieu 8  @("CT", 593, "07-04-2024", "LONG XUYÊN", "Trang thiết bị Y Tế", 40000)
Write-Row $wsChiTietChiTieu 9  @("CT", 597, "07-06-2024", "LONG XUYÊN", "Ứng Lương", 320000)
Write-Row $wsChiTietChiTieu 10 @("CT", 598, "07-06-2024", "LONG XUYÊN", "Chi Phí Sinh Hoạt Tại Cơ Sở", 365000)

# 3b. "DOANH SỐ CÁ NHÂN" (previously mis-spelled "DAONH SỐ CÁ NHÂN")
$wsDoanhSo = $wb.Worksheets.Add($null, $wsChiTietChiTieu)
$wsDoanhSo.Name = "DOANH SỐ CÁ NHÂN"
Set-StandardMargins $wsDoanhSo

Write-Row $wsDoanhSo 1 @("Mã nhân viên", "Doanh số sale chính", "Doanh số upsale", "Doanh số đơn 1 bác sĩ", "Doanh số đơn 2 bác sĩ", "Số lần phụ phẫu 1", "Công phụ phẫu 1", "Số lần phụ phẫu 2", "Công phụ phẫu 2", "Doanh số thu nợ")
Write-Row $wsDoanhSo 2 @("Nguyễn Phúc Nam", 3000000, 0, 0, 0, 0, 0, 0, 0, 800000)
Write-Row $wsDoanhSo 3 @("Đào Vương Anh", 0, 0, 0, 0, 1, 0, 0, 0, 0)
Write-Row $wsDoanhSo 4 @("Đặng Ngọc Mai", 0, 0, 4500000, 0, 0, 0, 0, 0, 0)
Write-Row $wsDoanhSo 5 @("Tổng", 3000000, 0, 4500000, 0, 1, 0, 0, 0, 800000)

# 3c. "CHI TIÊU TỔNG HỢP" (previously "CHI TIÊU")
$wsChiTieuTongHop = $wb.Worksheets.Add($null, $wsDoanhSo)
$wsChiTieuTongHop.Name = "CHI TIÊU TỔNG HỢP"
Set-StandardMargins $wsChiTieuTongHop

Write-Row $wsChiTieuTongHop 1 @("Phân loại", "Lượng chi")
Write-Row $wsChiTieuTongHop 2 @("Chi Phí CTV", 3080000)
Write-Row $wsChiTieuTongHop 3 @("Chi Phí Sinh Hoạt Tại Cơ Sở", 895000)
Write-Row $wsChiTieuTongHop 4 @("Trang thiết bị Y Tế", 40000)
Write-Row $wsChiTieuTongHop 5 @("Ứng Lương", 5320000)
Write-Row $wsChiTieuTongHop 6 @("Blank", 0)
Write-Row $wsChiTieuTongHop 7 @("Tổng cộng", 9335000)

# 3d. "LŨY KẾ NGÀY" (unchanged content, recreated to keep sheetId sequential)
$wsLuyKe = $wb.Worksheets.Add($null, $wsChiTieuTongHop)
$wsLuyKe.Name = "LŨY KẾ NGÀY"
Set-StandardMargins $wsLuyKe

Write-Row $wsLuyKe 1 @("Ngày", "Đơn giá", "Thanh toán lần đầu", "Số lượng đơn", "Thu nợ", "Lượng chi")
Write-Row $wsLuyKe 2 @("07-01-2024", 3000000, 3000000, 1, 0, 6800000)
Write-Row $wsLuyKe 3 @("07-02-2024", 1500000, 1500000, 1, 800000, 530000)
Write-Row $wsLuyKe 4 @("07-03-2024", 0, 0, 0, 0, 1280000)
Write-Row $wsLuyKe 5 @("07-04-2024", 0, 0, 0, 0, 40000)
Write-Row $wsLuyKe 6 @("07-06-2024", 0, 0, 0, 0, 685000)
Write-Row $wsLuyKe 7 @("Tổng", 4500000, 4500000, 2, 800000, 9335000)

# ---------------------------------------------------------------------------
# 4. Restore the originally active tab (first sheet).
# ---------------------------------------------------------------------------
$wsDoanhThu.Activate()
